$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ "E"=3; "G"=23.061728; "H"=69.18518399999999; "I"=0.6130144106248721; "J"=0.6902769593117909; "K"=3; "M"=8.206141000000001; "N"=24.618423; "O"=0.1726874201070044; "P"=0.185725470694853; "Q"=189.247791671648; "R"=1703.230125044832; "S"=0.105859877059225; "T"=0.1282020131779943 }
  3 = @{ "E"=3; "G"=23.061728; "H"=69.18518399999999; "I"=0.6130144106248721; "J"=0.6902769593117909; "K"=3; "M"=19.840792; "N"=59.522376; "O"=0.4175233137426826; "P"=0.4490466874939968; "Q"=457.562948408576; "R"=4118.066535677183; "S"=0.2559478080961141; "T"=0.3099665820323881 }
  4 = @{ "E"=3; "G"=23.061728; "H"=69.18518399999999; "I"=0.6130144106248721; "J"=0.6902769593117909; "K"=3; "M"=4.442508333333334; "N"=13.327525; "O"=0.09348673181306551; "P"=0.100545061469714; "Q"=102.4519188210667; "R"=922.0672693896; "S"=0.05730871380363182; "T"=0.06940393930513129 }
  5 = @{ "E"=3; "G"=23.061728; "H"=69.18518399999999; "I"=0.6130144106248721; "J"=0.6902769593117909; "K"=3; "M"=5.022911666666666; "N"=15.068735; "O"=0.1057005548822571; "P"=0.1136810388159715; "Q"=115.8370226246933; "R"=1042.53320362224; "S"=0.06479596335386877; "T"=0.07847140180529449 }
  6 = @{ "E"=3; "G"=23.061728; "H"=69.18518399999999; "I"=0.6130144106248721; "J"=0.6902769593117909; "K"=2; "M"=10.0078485; "N"=20.015697; "O"=0.2106019794549903; "P"=0.1510017415254648; "Q"=230.798279972208; "R"=1384.789679833248; "S"=0.1291020483120323; "T"=0.1042330229909828 }
  7 = @{ "E"=3; "G"=1.926013333333334; "H"=5.778040000000001; "I"=0.05119624723650278; "J"=0.05764887294340218; "K"=3; "M"=8.206141000000001; "N"=24.618423; "O"=0.1726874201070044; "P"=0.185725470694853; "Q"=15.80513698121334; "R"=142.24623283092; "S"=0.008840947854432017; "T"=0.01070686406244115 }
  8 = @{ "E"=3; "G"=1.926013333333334; "H"=5.778040000000001; "I"=0.05119624723650278; "J"=0.05764887294340218; "K"=3; "M"=19.840792; "N"=59.522376; "O"=0.4175233137426826; "P"=0.4490466874939968; "Q"=38.21362993589334; "R"=343.9226694230401; "S"=0.0213756267973743; "T"=0.02588703543299704 }
  9 = @{ "E"=3; "G"=1.926013333333334; "H"=5.778040000000001; "I"=0.05119624723650278; "J"=0.05764887294340218; "K"=3; "M"=4.442508333333334; "N"=13.327525; "O"=0.09348673181306551; "P"=0.100545061469714; "Q"=8.556330283444446; "R"=77.00697255100002; "S"=0.004786169835234331; "T"=0.005796309473754105 }
  10 = @{ "E"=3; "G"=1.926013333333334; "H"=5.778040000000001; "I"=0.05119624723650278; "J"=0.05764887294340218; "K"=3; "M"=5.022911666666666; "N"=15.068735; "O"=0.1057005548822571; "P"=0.1136810388159715; "Q"=9.674194842155556; "R"=87.06775357940002; "S"=0.005411471740787564; "T"=0.006553583762775914 }
  11 = @{ "E"=3; "G"=1.926013333333334; "H"=5.778040000000001; "I"=0.05119624723650278; "J"=0.05764887294340218; "K"=2; "M"=10.0078485; "N"=20.015697; "O"=0.2106019794549903; "P"=0.1510017415254648; "Q"=19.27524964898; "R"=115.65149789388; "S"=0.01078203100867456; "T"=0.008705080211433975 }
  12 = @{ "E"=2; "G"=12.632464; "H"=25.264928; "I"=0.3357893421386252; "J"=0.2520741677448068; "K"=3; "M"=8.206141000000001; "N"=24.618423; "O"=0.1726874201070044; "P"=0.185725470694853; "Q"=103.663780761424; "R"=621.982684568544; "S"=0.0579865951933474; "T"=0.04681659345441758 }
  13 = @{ "E"=2; "G"=12.632464; "H"=25.264928; "I"=0.3357893421386252; "J"=0.2520741677448068; "K"=3; "M"=19.840792; "N"=59.522376; "O"=0.4175233137426826; "P"=0.4490466874939968; "Q"=250.638090671488; "R"=1503.828544028928; "S"=0.1401998788491942; "T"=0.1131930700286116 }
  14 = @{ "E"=2; "G"=12.632464; "H"=25.264928; "I"=0.3357893421386252; "J"=0.2520741677448068; "K"=3; "M"=4.442508333333334; "N"=13.327525; "O"=0.09348673181306551; "P"=0.100545061469714; "Q"=56.11982659053333; "R"=336.7189595432; "S"=0.03139184817419935; "T"=0.0253448126908286 }
  15 = @{ "E"=2; "G"=12.632464; "H"=25.264928; "I"=0.3357893421386252; "J"=0.2520741677448068; "K"=3; "M"=5.022911666666666; "N"=15.068735; "O"=0.1057005548822571; "P"=0.1136810388159715; "Q"=63.45175080434666; "R"=380.71050482608; "S"=0.03549311978760076; "T"=0.0286560532479011 }
  16 = @{ "E"=2; "G"=12.632464; "H"=25.264928; "I"=0.3357893421386252; "J"=0.2520741677448068; "K"=2; "M"=10.0078485; "N"=20.015697; "O"=0.2106019794549903; "P"=0.1510017415254648; "Q"=126.423785893704; "R"=505.6951435748159; "S"=0.07071790013428345; "T"=0.03806363832304797 }
}

foreach ($row in $data.Keys) {
  foreach ($col in $data[$row].Keys) {
    $ws.Range("$col$row").Value = $data[$row][$col]
  }
}
